# Actualización automática 2025-11-03 08:30:05
#
# Monthly rollover update:
#  1) "VENTAS POR GRUPO": zero out the product-group figures that belonged
#     to the month that just closed, and refresh the "<n> de 58" coverage
#     counters for the columns whose totals were reset.
#  2) "VENTA MENSUAL": slide the four rolling month columns (C:F) one
#     column to the left - this month's figures move into last month's
#     slot, etc. - and open a brand-new (empty) column for the new month,
#     updating the month-name headers and column widths to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) VENTAS POR GRUPO - zero the cells for the group(s) that rolled off
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$zeroCells = @(
    "H2", "L2", "M2",
    "E3", "G3", "M3",
    "K4",
    "E11", "M11",
    "P16",
    "M17", "P17",
    "L30",
    "L32",
    "D35", "M35",
    "M40",
    "D46", "M46",
    "M47", "P47",
    "K53", "M53",
    "D57",
    "M58",
    "M59"
)

foreach ($addr in $zeroCells) {
    $wsGrupo.Range($addr).Value2 = 0
}

# Refresh the "<count> de 58" labels on the totals row for every column
# whose count dropped to zero along with the amounts above.
$countCells = @{
    "D60" = "0 de 58"
    "E60" = "0 de 58"
    "H60" = "0 de 58"
    "K60" = "0 de 58"
    "L60" = "0 de 58"
    "M60" = "0 de 58"
    "P60" = "0 de 58"
}

foreach ($addr in $countCells.Keys) {
    $wsGrupo.Range($addr).Value2 = $countCells[$addr]
}

# ---------------------------------------------------------------------
# 2) VENTA MENSUAL - roll the C:F month window one column to the left
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# -- month headers (row 1): C<-D, D<-E, E<-F, F<-new month name
$hdr = $wsMensual.Range("C1:F1").Value2
$wsMensual.Range("C1").Value2 = $hdr[1,2]
$wsMensual.Range("D1").Value2 = $hdr[1,3]
$wsMensual.Range("E1").Value2 = $hdr[1,4]
$wsMensual.Range("F1").Value2 = "noviembre"

# -- data rows 2:60 (includes the totals row 60): same left-shift, the
#    freshly opened column starts at 0
$dataRange = $wsMensual.Range("C2:F60")
$data = $dataRange.Value2
$nrows = $data.GetLength(0)
$shifted = New-Object 'object[,]' $nrows, 4
for ($i = 1; $i -le $nrows; $i++) {
    $shifted[$i - 1, 0] = $data[$i, 2]
    $shifted[$i - 1, 1] = $data[$i, 3]
    $shifted[$i - 1, 2] = $data[$i, 4]
    $shifted[$i - 1, 3] = 0
}
$dataRange.Value2 = $shifted

# -- column widths follow the same left-shift; the vacated column gets a
#    fresh width of 15.
# Note: Excel's ColumnWidth getter already reports the raw OOXML <col
# width> minus a ~0.8333 char padding, and the setter re-adds that same
# padding - so piping a getter value straight into another column's
# setter reproduces the exact original raw width (the padding cancels
# out). Only a brand-new literal target width needs the manual
# compensation below.
$w3 = $wsMensual.Columns.Item(3).ColumnWidth
$w4 = $wsMensual.Columns.Item(4).ColumnWidth
$w5 = $wsMensual.Columns.Item(5).ColumnWidth
$w6 = $wsMensual.Columns.Item(6).ColumnWidth
$pad = 0.8333333333333334

$wsMensual.Columns.Item(3).ColumnWidth = $w4
$wsMensual.Columns.Item(4).ColumnWidth = $w5
$wsMensual.Columns.Item(5).ColumnWidth = $w6
$wsMensual.Columns.Item(6).ColumnWidth = (15 - $pad)
